# 5.17 Add History Image & Action
#
# Adds two new trailing columns ("History" / "HistoryAction") to the
# header row of Sheet1, widens the new HistoryAction column, and moves
# the active selection to the newly added cell (R4), matching the
# author's "Add History Image & Action" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells: Q1 = "History", R1 = "HistoryAction"
$ws.Range("Q1").Value = "History"
$ws.Range("R1").Value = "HistoryAction"

# Give the new HistoryAction column (R / column 18) a bit more width,
# same as the author's custom width of 14 characters.
$ws.Range("R1").ColumnWidth = 13.29

# Move/leave the selection on the newly added cell.
$ws.Range("R4").Select()
